# Automatische test-sync: 2025-07-29 22:02:50
# Adds the newest mail-log entry (Testmail #16) to the "Logs" sheet and
# refreshes the "Dashboard" summary counts / ordering accordingly.

$wb = $excel.ActiveWorkbook

# --- 1. Append the new log row (row 18) on the "Logs" sheet ---------------
$ws = $wb.Worksheets.Item("Logs")

$ws.Range("A18").Value = "Wil je dit even doorsturen?"
$ws.Range("B18").Value = "mailmind.test@zohomail.eu"
$ws.Range("C18").Value = "Testmail #16: Wil je dit even doorsturen?"
$ws.Range("D18").Value = "Overig"
$ws.Range("F18").Value = "2025-07-29 22:02:06"
$ws.Range("G18").Value = "Nee"
$ws.Range("H18").Value = "Ja"
$ws.Range("I18").Value = "Nee"
$ws.Range("J18").Value = "Nee"

# --- 2. Extend the conditional-formatting ranges to cover the new row -----
$dRules = $ws.Range("D2:D17").FormatConditions
for ($i = 1; $i -le $dRules.Count(); $i++) {
    $dRules.Item($i).ModifyAppliesToRange($ws.Range("D2:D18"))
}

$gRules = $ws.Range("G2:G17").FormatConditions
for ($i = 1; $i -le $gRules.Count(); $i++) {
    $gRules.Item($i).ModifyAppliesToRange($ws.Range("G2:G18"))
}

$hRules = $ws.Range("H2:H17").FormatConditions
for ($i = 1; $i -le $hRules.Count(); $i++) {
    $hRules.Item($i).ModifyAppliesToRange($ws.Range("H2:H18"))
}

$iRules = $ws.Range("I2:I17").FormatConditions
for ($i = 1; $i -le $iRules.Count(); $i++) {
    $iRules.Item($i).ModifyAppliesToRange($ws.Range("I2:I18"))
}

$jRules = $ws.Range("J2:J17").FormatConditions
for ($i = 1; $i -le $jRules.Count(); $i++) {
    $jRules.Item($i).ModifyAppliesToRange($ws.Range("J2:J18"))
}

# --- 3. Refresh the "Dashboard" summary table ------------------------------
# The new entry's category ("Overig") now ties "Productinformatie" at 5,
# and on a tie the categories are listed alphabetically, so the two rows
# swap places while "Productinformatie" gains the tied count of 5.
$ws2 = $wb.Worksheets.Item("Dashboard")

$ws2.Range("A2").Value = "Overig"
$ws2.Range("B2").Value = 5
$ws2.Range("A3").Value = "Productinformatie"
$ws2.Range("B3").Value = 5
